$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 386 (shifts existing rows 386-403 down to 389-406)
$ws.Rows.Item(386).Insert()
$ws.Rows.Item(386).Insert()
$ws.Rows.Item(386).Insert()

# Row 386 - new "Especial" record
$ws.Range("A386").Value = 2
$ws.Range("B386").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C386").Value = "Coquimbo"
$ws.Range("D386").Value = 44706
$ws.Range("E386").Value = 4
$ws.Range("F386").Value = "Fruta"
$ws.Range("G386").Value = 100101
$ws.Range("H386").Value = "Berries"
$ws.Range("I386").Value = 100112025
$ws.Range("J386").Value = "Frutilla"
$ws.Range("K386").Value = "Sin especificar"
$ws.Range("L386").Value = "Especial"
$ws.Range("M386").Value = 500
$ws.Range("N386").Value = 18000
$ws.Range("O386").Value = 19000
$ws.Range("P386").Value = 18500
$ws.Range("Q386").Value = "`$/bandeja 7 kilos"
$ws.Range("R386").Value = "Provincia de Melipilla"
$ws.Range("S386").Value = 2643
$ws.Range("T386").Value = 7

# Row 387 - new "Primera" record
$ws.Range("A387").Value = 2
$ws.Range("B387").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C387").Value = "Coquimbo"
$ws.Range("D387").Value = 44706
$ws.Range("E387").Value = 4
$ws.Range("F387").Value = "Fruta"
$ws.Range("G387").Value = 100101
$ws.Range("H387").Value = "Berries"
$ws.Range("I387").Value = 100112025
$ws.Range("J387").Value = "Frutilla"
$ws.Range("K387").Value = "Sin especificar"
$ws.Range("L387").Value = "Primera"
$ws.Range("M387").Value = 500
$ws.Range("N387").Value = 16000
$ws.Range("O387").Value = 17000
$ws.Range("P387").Value = 16500
$ws.Range("Q387").Value = "`$/bandeja 7 kilos"
$ws.Range("R387").Value = "Provincia de Melipilla"
$ws.Range("S387").Value = 2357
$ws.Range("T387").Value = 7

# Row 388 - new "Segunda" record
$ws.Range("A388").Value = 2
$ws.Range("B388").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C388").Value = "Coquimbo"
$ws.Range("D388").Value = 44706
$ws.Range("E388").Value = 4
$ws.Range("F388").Value = "Fruta"
$ws.Range("G388").Value = 100101
$ws.Range("H388").Value = "Berries"
$ws.Range("I388").Value = 100112025
$ws.Range("J388").Value = "Frutilla"
$ws.Range("K388").Value = "Sin especificar"
$ws.Range("L388").Value = "Segunda"
$ws.Range("M388").Value = 400
$ws.Range("N388").Value = 12000
$ws.Range("O388").Value = 13000
$ws.Range("P388").Value = 12500
$ws.Range("Q388").Value = "`$/bandeja 7 kilos"
$ws.Range("R388").Value = "Provincia de Melipilla"
$ws.Range("S388").Value = 1786
$ws.Range("T388").Value = 7
